$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) Insert 4 new rows before row 15, copying the formatting of row 13 so the
#    borders/styles used across rows 10-14 continue into the new rows 15-18.
# ---------------------------------------------------------------------------
$ws.Range("B15:L18").EntireRow.Insert()
$ws.Range("B13:L13").Copy()
$ws.Range("B15:L18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row heights for the newly created / shifted rows.
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6
$ws.Rows.Item(18).RowHeight = 16.2

# Row 14 becomes a "middle" row of the table (it used to be the last row
# before the thick bottom border) so its L cell switches to the wrapped
# style used by the other middle rows (style 60).
$ws.Range("L14").WrapText = $true

# Row 18 is now the last row before the thick closing border row (row 19);
# match the "no-wrap" look used by the old last row (style 44) by turning
# wrap text off for L18.
$ws.Range("L18").WrapText = $false

# ---------------------------------------------------------------------------
# 2) Update cell text content.
# ---------------------------------------------------------------------------
$ws.Range("D7").Value = "Bound Flasher"

$ws.Range("C10").Value = "addOp check"
$ws.Range("E10").Value = "addOp check in normal cases"
$ws.Range("G10").Value = "Check addition operator between 2 positive parameters"
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = ""
$ws.Range("K10").Value = ""

$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("G11").Value = "Check addition operator between 2 negative parameters"
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""

$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("G12").Value = "Check addition operator between a positive parameter and a negative one"
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""

$ws.Range("D13").Value = 2
$ws.Range("E13").Value = "addOp check in overflow cases"
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = "Check addition operator between 2 positive parameters"

$ws.Range("D14").Value = 3
$ws.Range("E14").Value = "addOp check in underflow cases"
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = "Check addition operator between 2 negative parameters"

$ws.Range("D15").Value = 4
$ws.Range("E15").Value = "addOp check in special cases"
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = "Add a number with zero"

$ws.Range("D16").Value = 7
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = "Add a number with its multiplicative inverse value"

$ws.Range("D17").Value = 8

$ws.Range("D18").Value = 9
$ws.Range("F18").Value = 5

# ---------------------------------------------------------------------------
# 3) Selection / view tweaks.
# ---------------------------------------------------------------------------
$ws.Range("G15").Select()

# ---------------------------------------------------------------------------
# 4) Conditional formatting range grows with the new rows (Excel does not
#    auto-extend FormatConditions ranges when rows are inserted inside them).
# ---------------------------------------------------------------------------
$fcs = $ws.Range("H10:J14").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("H10:J18"))
}
